# Update the "2024" sheet: insert a new transaction entry at the top of the
# "Others" category list (row 45), pushing the existing entries down by one
# row. This also naturally shifts the "Broadband" label in column A (row 178)
# down to row 179, matching the new sheet dimension A1:Y179.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row before row 45, shifting rows 45:178 down to 46:179.
$ws.Rows.Item(45).Insert()

# Populate the new row with the latest transaction entry.
$ws.Cells.Item(45, 18).Value = "bal axis"
$ws.Cells.Item(45, 19).Value = "2024-09-20 07:03:45"
